# Applies the "1.2.4 -> 1.2.5" version bump and accompanying text fixes
# to the "UC009 - Prestar Contas" test-suite workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# 1) Version: "1.0" -> "1.2.5"
$ws.Range("D2").Value = "1.2.5"

# 2) Precondition text: fix "usuario" -> "usuário" and add trailing period.
#    This text is repeated (as the same shared string) for every test case.
$preconditionCells = @("B8", "B18", "B28", "B37", "B46", "B55", "B64", "B72", "B80", "B87")
foreach ($cellRef in $preconditionCells) {
    $ws.Range($cellRef).Value = "O usuário devidamente autenticado e na tela inicial do sistema."
}

# 3) TC2 step: fix typo "histório" -> "histórico"
$ws.Range("B22").Value = "Chefe Verifica o histórico da tramitação da prestação de contas."

# 4) Swap TC3 / TC4 step-3 content (delete-comprovante and view-comprovante steps were reordered)
$ws.Range("B32").Value = "Chefe Clica em excluir comprovante."
$ws.Range("D32").Value = "SYSTEM Exclui o comprovante."
$ws.Range("B41").Value = "Chefe Clica em visualizar comprovante."
$ws.Range("D41").Value = "SYSTEM Exibe modal com o comprovante."

# 5) TC5 expected result: add trailing period
$ws.Range("D50").Value = "SYSTEM Apresenta a tela de Detalhar Diárias."

# 6) TC7 expected result: fix "Permite não permite" -> "Não permite"
$ws.Range("D67").Value = "SYSTEM Identifica que a prestação de contas indicada pelo usuário não está em nenhum desses dois estados: a) NÃO REALIZADA e b) DEVOLVIDA; Não permite um novo envio ou alterações na prestação (exclusão de documentos)."
